$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert one new row before row 13 -------------------------------------
# This shifts the former rows 13-23 down to 14-24, carrying their styles
# and row heights with them, matching the target layout. It also shifts the
# "answer" column (B/C) content relative to the "label" column (A) content
# by one slot starting at row 10, which is corrected below cell by cell.
$ws.Rows.Item(13).Insert()

# The newly inserted row 13 picks up an empty A13 cell (style copied from the
# row above); the target layout has no A13 at all (only B13/C13 hold data),
# so drop it.
$ws.Range("A13").Clear()

# --- New text values --------------------------------------------------------
$objetivosText = 'Introduzir as técnicas de fabricação de dispositivos e circuitos integrados em microeletrônica. Apresenta os princípios, técnicas, equipamentos e softwares utilizados na simulação e fabricação de dispositivos em silício e arseneto de gálio de uma maneira global e genérica.'

$docentesText = '5840535 - Messias Borges Silva'

$programaResumidoText = 'Processos de fabricação típicos e de montagem de sistemas eletroeletronicos'

$programaText = 'Técnicas de obtenção de silício cristalino e arseneto de gálio. Processamento de lâminas de silício e GaAs. Rede e Defeitos cristalinos. Oxidação. Dopagem (difusão e implantação iônica). Litografia: (Fabricação de Fotomáscaras; Gerador de Padrões). Epitaxia. Deposição em Fase Vapor (CVD. PECVD e LPCVD). Decapagem (úmida e seca). Equipamentos e técnicas de metalização (Evaporação térmica, feixe iônico, bombardeamento catódico ("sputtering"), Caracterização de etapas de processo de fabricação. Simulação de processos de fabricação. Montagem de Sistemas Eletro Eletrônicos'

$metodoText = 'Aulas expositivas auxiliadas por transparências e métodos multimídia (datashow + computador) São previstas três aulas práticas com visitas em ambientes de fabricação de circuitos integrados (salas limpas) e montagem de produtos eletroeletrônicos. Duas prvas e relatórios de trabalhos práticos e visitas.'

$criterioText = 'Média geral MG = [0,6 x (média aritmética das 2 provas feitas) + 0,2 x média dos trabalhos práticos + 0,2 x (média aritmética dos relatórios de visitas)] >= 5,0. Os testes serão dados nos 10 minutos finais de N aulas escolhidas aleatoriamente. A prova substitutiva substitui a prova em que o aluno faltou.'

$normaRecuperacaoText = '1 (uma) prova de recuperação.'

$bibliografiaText = '01 Stephen A. Campbell, The Science and Engineering of Microelectronic Fabrication, 2nd edition, Oxford University Press, 2001;02 S.M. Sze - VLSI Technology, McGraw-Hill, 1988;03 V. Baranauskas, org., Processos em Microeletrônica, SBV e SBMicro, 1990;04 R. A. Levy, Microelectronic Materials and Processes, Kluwer Academic Publ., 1989;[05] M. Madou, Fundamentals of Microfabrication, CRC press, 1997.'

$requisitosText = 'LOQ4205 -  Sistemas Produtivos II  (Requisito fraco)'

# Row 10 (Objetivos:): replace the stray "Docentes" text with the real objectives text.
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# Row 13 (new, blank row under "Docentes responsáveis:"): holds the teacher entry
# that used to sit on row 10. Copy B9:C9's cell format first (wrap-text body
# style) since the freshly inserted row has no format of its own there.
$ws.Range("B9:C9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = $docentesText
$ws.Range("C13").Value = $docentesText

# Row 14 (Programa resumido:): replace stray "Semestral" with the real short syllabus.
$ws.Range("B14").Value = $programaResumidoText
$ws.Range("C14").Value = $programaResumidoText

# Row 15 (Short syllabus:): clear the stray duplicated date that used to live here.
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# Row 16 (Programa:): fill in the full syllabus text (previously empty).
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# Row 19 (Método:): replace the stray teacher text with the real method text.
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20 (Critério:): replace the method text (that had slid here) with the
# real grading-criteria text.
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21 (Norma de recuperação:): replace the criteria text with the real
# recovery-exam text.
$ws.Range("B21").Value = $normaRecuperacaoText
$ws.Range("C21").Value = $normaRecuperacaoText

# Row 22 (Bibliografia:): replace the recovery-exam text with the real
# bibliography text.
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText

# Row 24 (prerequisite entry under "Requisitos:"): same text as before, but
# the trailing line break is dropped.
$ws.Range("B24").Value = $requisitosText
$ws.Range("C24").Value = $requisitosText
